# Actualización 11 de Mayo - Mañana
# Updates to the "Blancos" (J) / "Por_Blan" (K) columns for Castro Vasquez
# Julieta's "2ARHM" row on all three partial-exam sheets, plus corrected
# Promedio/Aprobados/Reprobados figures for Martínez López Miguel Ángel's
# 6ARHV rows on the "2o Parcial" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "1er Parcial": row 5 (Castro Vasquez Julieta / 2ARHM) ---
$ws1 = $wb.Worksheets.Item("1er Parcial")
$ws1.Cells.Item(5, 10).Value = 0   # J5 Blancos
$ws1.Cells.Item(5, 11).Value = 0   # K5 Por_Blan

# --- Sheet "2o Parcial" ---
$ws2 = $wb.Worksheets.Item("2o Parcial")

# Row 5 (Castro Vasquez Julieta / 2ARHM)
$ws2.Cells.Item(5, 9).Value = 6.5    # I5 Promedio
$ws2.Cells.Item(5, 10).Value = 0     # J5 Blancos
$ws2.Cells.Item(5, 11).Value = 0     # K5 Por_Blan

# Row 9 (Martínez López Miguel Ángel / 6ARHV / DETERMINA LA NÓMINA...)
$ws2.Cells.Item(9, 5).Value = 22          # E9 Aprobados
$ws2.Cells.Item(9, 6).Value = 6           # F9 Reprobados
$ws2.Cells.Item(9, 7).Value = 78.56999999999999   # G9 Por_Apro
$ws2.Cells.Item(9, 8).Value = 21.43       # H9 Por_Repro
$ws2.Cells.Item(9, 9).Value = 7.9         # I9 Promedio
$ws2.Cells.Item(9, 10).Value = 6          # J9 Blancos
$ws2.Cells.Item(9, 11).Value = 21.43      # K9 Por_Blan

# Row 10 (Martínez López Miguel Ángel / 6ARHV / DETERMINA REMUNERACIONES...)
$ws2.Cells.Item(10, 5).Value = 19     # E10 Aprobados
$ws2.Cells.Item(10, 6).Value = 9      # F10 Reprobados
$ws2.Cells.Item(10, 7).Value = 67.86  # G10 Por_Apro
$ws2.Cells.Item(10, 8).Value = 32.14  # H10 Por_Repro
$ws2.Cells.Item(10, 9).Value = 7.9    # I10 Promedio
$ws2.Cells.Item(10, 10).Value = 8     # J10 Blancos
$ws2.Cells.Item(10, 11).Value = 28.57 # K10 Por_Blan

# --- Sheet "3er Parcial": row 5 (Castro Vasquez Julieta / 2ARHM) ---
$ws3 = $wb.Worksheets.Item("3er Parcial")
$ws3.Cells.Item(5, 10).Value = 0   # J5 Blancos
$ws3.Cells.Item(5, 11).Value = 0   # K5 Por_Blan
